{"js": "// Replace Product, IT, and Finance templates with correct industry-specific content.\n// Applies a series of exact text replacements plus removal of a page-break run.\n\nasync function replaceAll(body, searchText, replacement) {\n  const results = body.search(searchText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replacement, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\nconst body = context.document.body;\n\n// 1. Title block: \"ARTIFICIAL INTELLIGENCE AND MACHINE LEARNING\" -> \"PRODUCT DEVELOPMENT\"\nawait replaceAll(body, \"ARTIFICIAL INTELLIGENCE AND MACHINE LEARNING\", \"PRODUCT DEVELOPMENT\");\n\n// 2. \"Product and Machine Learning\" -> \"Product Development and Product Innovation\" (4 occurrences)\nawait replaceAll(body, \"Product and Machine Learning\", \"Product Development and Product Innovation\");\n\n// 3. \"Product Implementation\" -> \"Product Development Implementation\" (5 occurrences)\nawait replaceAll(body, \"Product Implementation\", \"Product Development Implementation\");\n\n// 4. Stakeholder / role renames\nawait replaceAll(body, \"ML Engineers\", \"Product Engineers\");\nawait replaceAll(body, \"Compliance Officers\", \"Quality Assurance Managers\");\n\n// 5. Technology component renames\nawait replaceAll(body, \"MLflow\", \"Productflow\");\nawait replaceAll(body, \"Cloud ML platforms\", \"Cloud Product platforms\");\n\n// 6. \"Financial\" -> \"Product\" (section labels, 4 occurrences)\nawait replaceAll(body, \"Financial Justification:\", \"Product Justification:\");\nawait replaceAll(body, \"Financial: Budget overruns, cost escalation, ROI delays\", \"Product: Budget overruns, cost escalation, ROI delays\");\nawait replaceAll(body, \"ML Platform Licensing: $840,000\", \"Product Platform Licensing: $840,000\");\nawait replaceAll(body, \"Financial: Break-even within 30 months, 250%+ ROI within 3 years\", \"Product: Break-even within 30 months, 250%+ ROI within 3 years\");\nawait replaceAll(body, \"Financial review and budget allocation approval\", \"Product review and budget allocation approval\");\n\n// 7. Remove the standalone page-break paragraph's run (keep paragraph, empty run -> no more page break).\nconst paras = body.paragraphs;\nparas.load(\"items/text\");\nawait context.sync();\nfor (let i = 0; i < paras.items.length; i++) {\n  if (paras.items[i].text === \"\\f\") {\n    paras.items[i].clear();\n  }\n}\nawait context.sync();\n", "ps1": "# Replace Product, IT, and Finance templates with correct industry-specific content.\n$d = $word.ActiveDocument\n\nfunction Replace-All {\n    param(\n        [string]$FindText,\n        [string]$ReplaceText\n    )\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Execute($FindText, $false, $false, $false, $false, $false, $true, 1, $false, $ReplaceText, 2)\n}\n\n# 1. Title block\nReplace-All \"ARTIFICIAL INTELLIGENCE AND MACHINE LEARNING\" \"PRODUCT DEVELOPMENT\"\n\n# 2. \"Product and Machine Learning\" -> \"Product Development and Product Innovation\" (4 occurrences)\nReplace-All \"Product and Machine Learning\" \"Product Development and Product Innovation\"\n\n# 3. \"Product Implementation\" -> \"Product Development Implementation\" (5 occurrences)\nReplace-All \"Product Implementation\" \"Product Development Implementation\"\n\n# 4. Stakeholder / role renames\nReplace-All \"ML Engineers\" \"Product Engineers\"\nReplace-All \"Compliance Officers\" \"Quality Assurance Managers\"\n\n# 5. Technology component renames\nReplace-All \"MLflow\" \"Productflow\"\nReplace-All \"Cloud ML platforms\" \"Cloud Product platforms\"\n\n# 6. \"Financial\" -> \"Product\" (section labels, 4 occurrences)\nReplace-All \"Financial Justification:\" \"Product Justification:\"\nReplace-All \"Financial: Budget overruns, cost escalation, ROI delays\" \"Product: Budget overruns, cost escalation, ROI delays\"\nReplace-All \"ML Platform Licensing: `$840,000\" \"Product Platform Licensing: `$840,000\"\nReplace-All \"Financial: Break-even within 30 months, 250%+ ROI within 3 years\" \"Product: Break-even within 30 months, 250%+ ROI within 3 years\"\nReplace-All \"Financial review and budget allocation approval\" \"Product review and budget allocation approval\"\n\n# 7. Remove the standalone page-break paragraph's content (keep paragraph, drop the page break).\n$paras = $d.Paragraphs\nfor ($i = 1; $i -le $paras.Count; $i++) {\n    $p = $paras.Item($i)\n    $t = $p.Range.Text\n    if ($t.Length -eq 2 -and [int][char]$t[0] -eq 12) {\n        $p.Range.Text = \"\"\n    }\n}\n"}
